# q695 and fix qnnn/test_.py to fix pytest bug
# Adds 3 new LeetCode tracker rows (Minimum Number of Taps to Open to
# Water a Garden / Flood Fill / Max Area of Island) below the existing
# data and keeps column widths + conditional formatting consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- widen column C a touch (bestFit grew with the new dates) ---------
$ws.Columns.Item(3).ColumnWidth = 15.6640625

# --- row 38: Minimum Number of Taps to Open to Water a Garden (Hard) --
$ws.Range("A38").Value = 1326
$ws.Range("B4").Copy($ws.Range("B38"))
$ws.Range("B38").Value = "Minimum Number of Taps to Open to Water a Garden"
$ws.Range("C37").Copy($ws.Range("C38"))
$ws.Range("C38").Value = 44544.671527777777
$ws.Range("D38").Value = "Hard"

# --- row 39: Flood Fill (Easy) -----------------------------------------
$ws.Range("A39").Value = 733
$ws.Range("B7").Copy($ws.Range("B39"))
$ws.Range("B39").Value = "Flood Fill"
$ws.Range("C37").Copy($ws.Range("C39"))
$ws.Range("C39").Value = 44545.451388888891
$ws.Range("D39").Value = "Easy"

# --- row 40: Max Area of Island (Medium) --------------------------------
$ws.Range("A40").Value = 695
$ws.Range("B7").Copy($ws.Range("B40"))
$ws.Range("B40").Value = "Max Area of Island"
$ws.Range("C37").Copy($ws.Range("C40"))
$ws.Range("C40").Value = 44545.58037152778
$ws.Range("D40").Value = "Medium"

# --- conditional formatting: keep the red/yellow/green Hard/Medium/Easy
# coloring going for the freshly-added rows the same way it was applied
# to each of the D31:D37 rows already in the sheet -----------------------
$rngs = @($ws.Range("D38"), $ws.Range("D39"), $ws.Range("D40"))
foreach ($r in $rngs) {
    $fcHard = $r.FormatConditions.Add(1, 3, '"Hard"')
    $fcHard.Interior.Color = 255
    $fcMedium = $r.FormatConditions.Add(1, 3, '"Medium"')
    $fcMedium.Interior.Color = 65535
    $fcEasy = $r.FormatConditions.Add(1, 3, '"Easy"')
    $fcEasy.Interior.Color = 5287936
}

# --- keep the view pointed at the new rows, like the saved workbook ----
$ws.Range("A11").Select()
$ws.Range("D40").Select()
